$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Step02A" (physical sheet1.xml): keep existing File/comment
# header and first data row, add a new data row for the second
# preprocessed-dataset scenario.
# ------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Step02A")
$wsA.Range("A3").Value = "PreprocessedDatasetScenario2.mat"
$wsA.Range("B3").Value = "cropped to only face and center images"

# ------------------------------------------------------------------
# Sheet "Step02B" (physical sheet2.xml): the old wide table
# (File / flatDataFile / numTrainingSamples / numTestSamples /
# Comment) is replaced by a narrower one
# (File / preprocessedDatasetFile / Comment). Remove the two
# now-unused columns (old D = numTestSamples, old E = Comment)
# so everything to the right shifts left by two columns, then
# update headers/values and drop the now-obsolete extra rows.
# ------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Step02B")

# Drop the "numTestSamples" / "Comment" columns entirely (this
# shifts the old F:M style-only columns left to D:K).
$wsB.Range("D1:E1").EntireColumn.Delete()

# Header row
$wsB.Range("A1").Value = "File"
$wsB.Range("B1").Value = "preprocessedDatasetFile"
$wsB.Range("C1").Value = "Comment"

# Row 2
$wsB.Range("A2").Value = "TrainingAndTestDataScenario1.mat"
$wsB.Range("B2").Value = "PreprocessedDatasetScenario1.mat"
$wsB.Range("C2").Value = "Use subject 14 and 15 as test"

# Row 3
$wsB.Range("A3").Value = "TrainingAndTestDataScenario2.mat"
$wsB.Range("B3").Value = "PreprocessedDatasetScenario2.mat"
$wsB.Range("C3").Value = "^"

# Remove the now-obsolete rows 4-6 of the old wide table.
$wsB.Range("A4:A6").EntireRow.Delete()

# ------------------------------------------------------------------
# Selection / active sheet bookkeeping to match the saved view state.
# ------------------------------------------------------------------
$null = $wsA.Range("A4").Select()
$null = $wsB.Range("B7").Select()
